$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.848.59"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "2.219.49"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "292.81"
$ws.Range("E5").Value = "  -1.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "86.15"
$ws.Range("E6").Value = "  +5.47%  "
$ws.Range("E7").Value = "  +1.17%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.473"
$ws.Range("E9").Value = "  +1.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.75"
$ws.Range("E10").Value = "  +5.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0786"
$ws.Range("E11").Value = "  +2.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.39"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("E13").Value = "  +1.63%  "
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("D15").Value = "2.562.05"
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").Value = "2.201.78"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.733"
$ws.Range("E18").Value = "  +3.37%  "
$ws.Range("D19").Value = "39.816.40"
$ws.Range("E19").Value = "  +1.97%  "
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.09"
$ws.Range("E21").Value = "  +7.98%  "
$ws.Range("E22").Value = "  +1.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.60"
$ws.Range("E23").Value = "  +1.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.38"
$ws.Range("E24").Value = "  +4.87%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  +2.89%  "
$ws.Range("E27").Value = "  +2.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.79"
$ws.Range("E28").Value = "  +1.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.12"
$ws.Range("E29").Value = "  -2.34%  "
$ws.Range("E30").Value = "  +1.99%  "
$ws.Range("E31").Value = "  +4.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "151.77"
$ws.Range("E32").Value = "  +1.58%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  +2.76%  "
$ws.Range("E35").Value = "  +4.03%  "
$ws.Range("E36").Value = "  +1.83%  "
$ws.Range("E37").Value = "  +7.21%  "
$ws.Range("E38").Value = "  +2.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.88"
$ws.Range("E39").Value = "  +4.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0993"
$ws.Range("E40").Value = "  +3.30%  "
$ws.Range("E41").Value = "  +4.35%  "
$ws.Range("E42").Value = "  +5.33%  "
$ws.Range("D43").Value = "2.065.43"
$ws.Range("E43").Value = "  +9.17%  "
$ws.Range("E44").Value = "  +3.98%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.09"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.92"
$ws.Range("E46").Value = "  +10.83%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.67"
$ws.Range("E47").Value = "  +10.15%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").Value = "2.432.47"
$ws.Range("E49").Value = "  +1.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.50"
$ws.Range("E50").Value = "  +0.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "89.05"
$ws.Range("E51").Value = "  +2.60%  "
